$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = [double]"0.55291933729290854"
$ws.Range("C5").Value = [double]"3.7201655762508425E-2"
$ws.Range("D5").Value = [double]"4.2080889838379265E-2"
$ws.Range("B7").Value = [double]"5.8790176509471062"
$ws.Range("C7").Value = [double]"0.13985954595569319"
$ws.Range("D7").Value = [double]"1.5753284961181173"
$ws.Range("B8").Value = [double]"0.55051201350063084"
$ws.Range("C8").Value = [double]"4.6957591701251286E-2"
$ws.Range("D8").Value = [double]"0.11878642186747311"
$ws.Range("B11").Value = [double]"1.6689695754745895"
$ws.Range("C11").Value = [double]"3.1380121511281181E-2"
$ws.Range("D11").Value = [double]"0.24742220209809906"
$ws.Range("B12").Value = [double]"4.4544770248008625"
$ws.Range("C12").Value = [double]"0.12022237610802837"
$ws.Range("D12").Value = [double]"1.7308597322619541"
$ws.Range("B14").Value = [double]"3.0698578776982726"
$ws.Range("C14").Value = [double]"7.8279827204200209E-2"
$ws.Range("D14").Value = [double]"0.71602876309792463"
$ws.Range("B17").Value = [double]"0.49573740669255484"
$ws.Range("C17").Value = [double]"2.9755271303232591E-3"
$ws.Range("D17").Value = [double]"0.1608022468984035"
$ws.Range("B18").Value = [double]"1.0568129784875988"
$ws.Range("C18").Value = [double]"3.8452801188592407E-2"
$ws.Range("D18").Value = [double]"1.0183238922906932"
$ws.Range("B19").Value = [double]"1.2211561218296874"
$ws.Range("C19").Value = [double]"8.2775109122183291E-2"
$ws.Range("D19").Value = [double]"0.88921937428561226"
$ws.Range("B20").Value = [double]"10.123514171923757"
$ws.Range("C20").Value = [double]"0.49869438826167795"
$ws.Range("D20").Value = [double]"2.0090421342222715"
$ws.Range("B21").Value = [double]"2.1916618435678838"
$ws.Range("C21").Value = [double]"0.10319873701097287"
$ws.Range("D21").Value = [double]"0.18750519783263128"
$ws.Range("B22").Value = [double]"1.1007247072784656"
$ws.Range("C22").Value = [double]"2.1255279380347446E-2"
$ws.Range("D22").Value = [double]"0.45222927090051473"
$ws.Range("B23").Value = [double]"2.3217323939712551"
$ws.Range("C23").Value = [double]"9.796895816614698E-2"
$ws.Range("D23").Value = [double]"1.2629865015927084"
$ws.Range("B24").Value = [double]"0.41773267501701572"
$ws.Range("C24").Value = [double]"3.2595782265507578E-2"
$ws.Range("D24").Value = [double]"5.3206994467569664E-2"
$ws.Range("B25").Value = [double]"2.3844277034814962"
$ws.Range("C25").Value = [double]"1.9867909290122648E-2"
$ws.Range("D25").Value = [double]"1.5680209108130214"
$ws.Range("B26").Value = [double]"1.9187149638260326"
$ws.Range("C26").Value = [double]"0.1059317666762583"
$ws.Range("D26").Value = [double]"0.48706006072605651"
$ws.Range("B27").Value = [double]"3.2394831666078354"
$ws.Range("C27").Value = [double]"0.16024184659131152"
$ws.Range("D27").Value = [double]"0.73531212416477065"
$ws.Range("B28").Value = [double]"2.1395774037265842"
$ws.Range("C28").Value = [double]"9.5705924763737599E-2"
$ws.Range("D28").Value = [double]"0.32430026545089363"
